$d = $word.ActiveDocument

$replacements = @(
    @{old="811×4="; new="814×6="},
    @{old="870×9="; new="118×2="},
    @{old="153×4="; new="843×2="},
    @{old="585×7="; new="809×6="},
    @{old="481×9="; new="819×4="},
    @{old="711×5="; new="367×6="},
    @{old="169×6="; new="358×6="},
    @{old="330×3="; new="448×2="},
    @{old="305×4="; new="600×4="},
    @{old="281×2="; new="104×7="},
    @{old="204×5="; new="710×8="},
    @{old="567×3="; new="950×2="},
    @{old="400×6="; new="692×9="},
    @{old="666×2="; new="719×8="},
    @{old="382×2="; new="432×8="},
    @{old="605×7="; new="738×8="},
    @{old="105×2="; new="390×7="},
    @{old="598×6="; new="378×4="},
    @{old="300×3="; new="896×9="},
    @{old="834×5="; new="559×3="},
    @{old="745×2="; new="762×4="},
    @{old="611×3="; new="669×7="},
    @{old="637×4="; new="162×8="},
    @{old="564×4="; new="548×3="},
    @{old="965×5="; new="971×3="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
